$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-03-25 Monday" "2024-03-26 Tuesday"

Replace-Text "862×9=7758" "514×6=3084"
Replace-Text "646×9=5814" "846×5=4230"
Replace-Text "756×6=4536" "169×8=1352"
Replace-Text "144×3=432" "123×6=738"
Replace-Text "610×9=5490" "202×3=606"

Replace-Text "725×8=5800" "544×3=1632"
Replace-Text "908×3=2724" "975×7=6825"
Replace-Text "572×9=5148" "411×2=822"
Replace-Text "236×6=1416" "474×8=3792"
Replace-Text "685×6=4110" "930×3=2790"

Replace-Text "978×3=2934" "939×4=3756"
Replace-Text "138×2=276" "114×8=912"
Replace-Text "197×4=788" "144×8=1152"
Replace-Text "596×8=4768" "570×6=3420"
Replace-Text "804×6=4824" "325×4=1300"

Replace-Text "448×6=2688" "788×8=6304"
Replace-Text "641×6=3846" "895×6=5370"
Replace-Text "534×6=3204" "130×3=390"
Replace-Text "144×6=864" "647×5=3235"
Replace-Text "519×4=2076" "340×9=3060"

Replace-Text "879×6=5274" "944×8=7552"
Replace-Text "475×2=950" "250×6=1500"
Replace-Text "258×2=516" "757×2=1514"
Replace-Text "389×7=2723" "556×7=3892"
Replace-Text "855×4=3420" "167×9=1503"
